# Minor text edit on one slide:
# "Once a loop statement has been parsed, we don't need to retain the
#  nonterminal symbols. ..." becomes "... retain the terminal symbols. ..."
# (the word "non" is removed from "nonterminal"), which also splits the
# single run that held the whole sentence into three runs at the edit
# boundaries.

$p = $ppt.ActivePresentation

# Locate the shape (on whichever slide it lives on) whose text contains the
# unique sentence we need to edit, rather than hard-coding slide/shape
# indices.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -like "*Once a loop statement has been parsed*") {
                $targetShape = $sh
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph containing the sentence that needs editing.
$targetPara = $null
for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
    $para = $tr.Paragraphs($pi, 1)
    if ($para.Text -like "*Once a loop statement has been parsed*") {
        $targetPara = $para
    }
}

# Replace "the nonterminal " with "the terminal " (i.e. drop "non"), using a
# Characters() sub-range so the edit lands exactly on that span and the
# surrounding, unmodified text keeps its original run/formatting.
$needle = "the nonterminal "
$startIdx = $targetPara.Text.IndexOf($needle)
$chars = $targetPara.Characters($startIdx + 1, $needle.Length)
$chars.Text = "the terminal "
